# [Kadastro App] Yeni kayit eklendi: 2953
# Append the new record row to both the "Kayitlar" (all records) sheet
# and the "Erdemli" (unit-filtered) sheet, one row below the current
# last row of each.

$wb = $excel.ActiveWorkbook

$newRecord = @("2953", "2025-09-09", "Erdemli", "1", "LİHKAB", "EMİNE ALANLI KIRCILI (K.Mühendisi), AYHAN KARADAYI (K.Teknisyeni)")

$targetSheets = @("Kayitlar", "Erdemli")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $newRow = $lastRow + 1

    for ($col = 1; $col -le $newRecord.Count; $col++) {
        $cell = $ws.Cells.Item($newRow, $col)
        # Prefix with a single quote so numeric-looking values ("2953", "1")
        # and the date-looking value ("2025-09-09") are entered as literal
        # text (quote-prefixed, General format) instead of being coerced
        # into numbers/dates - matching how the rest of the sheet stores
        # every column as text.
        $cell.Value = "'" + $newRecord[$col - 1]
    }
}
